$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8 (current row 8 "Basic information > Mobile No"
# test case and all following rows shift down by one).
$ws.Rows("8:8").Insert() | Out-Null

# Copy the formatting (styles) from the row that is now directly below
# (row 9, formerly row 8) onto the freshly inserted, blank row 8 so the
# new row reuses the existing cell styles instead of minting new ones.
$ws.Range("B9:G9").Copy() | Out-Null
$ws.Range("B8:G8").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Fill in the new test case row.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Buyer users page"
$ws.Range("C8").Value = "View users  page"
$ws.Range("D8").Value = "Basic information > Mobile phone"
$ws.Range("E8").Value = "If market is AU prefix should be set to +61. need to allow longer numbers (up to 20 chars)."
$ws.Range("F8").Value = "It gets displayed as expected"
$ws.Range("G8").Value = "Pass"
$ws.Rows("8:8").RowHeight = 45

# Renumber the "SL. No" column for the rows that shifted down (they kept
# their old sequence numbers after the insert) so numbering stays
# contiguous 8..17.
for ($r = 9; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Restore selection to just past the last row, like the source workbook.
$ws.Range("A19").Select() | Out-Null
